$wb = $excel.ActiveWorkbook

# zh-cn sheet: refresh the handback datetime for the first file's row (row 2)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-01-15 10:08:01"
$wsZh.Range("G2").Value = "2016-01-15 10:08:47"

# de-de sheet: refresh the handback datetime for the first file's row (row 2)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-01-15 10:08:12"
$wsDe.Range("G2").Value = "2016-01-15 10:09:05"
